$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("G2").Value = 63.91118233333333
    $ws.Range("H2").Value = 191.733547
    $ws.Range("I2").Value = 0.4067926910433548
    $ws.Range("J2").Value = 0.4067926910433549
    $ws.Range("M2").Value = 9.084137666666667
    $ws.Range("N2").Value = 27.252413
    $ws.Range("O2").Value = 0.2765376761551382
    $ws.Range("P2").Value = 0.2765376761551382
    $ws.Range("Q2").Value = 580.5779787554345
    $ws.Range("R2").Value = 5225.201808798911
    $ws.Range("S2").Value = 0.1124935054580244
    $ws.Range("T2").Value = 0.1124935054580244
    $ws.Range("G3").Value = 63.91118233333333
    $ws.Range("H3").Value = 191.733547
    $ws.Range("I3").Value = 0.4067926910433548
    $ws.Range("J3").Value = 0.4067926910433549
    $ws.Range("O3").Value = 0.3707916163717078
    $ws.Range("P3").Value = 0.3707916163717078
    $ws.Range("Q3").Value = 778.4597388884464
    $ws.Range("R3").Value = 7006.137649996017
    $ws.Range("S3").Value = 0.1508353194401623
    $ws.Range("T3").Value = 0.1508353194401623
    $ws.Range("G4").Value = 63.91118233333333
    $ws.Range("H4").Value = 191.733547
    $ws.Range("I4").Value = 0.4067926910433548
    $ws.Range("J4").Value = 0.4067926910433549
    $ws.Range("M4").Value = 11.58507333333333
    $ws.Range("N4").Value = 34.75522
    $ws.Range("O4").Value = 0.3526707074731541
    $ws.Range("P4").Value = 0.3526707074731541
    $ws.Range("Q4").Value = 740.4157341517044
    $ws.Range("R4").Value = 6663.741607365339
    $ws.Range("S4").Value = 0.1434638661451681
    $ws.Range("T4").Value = 0.1434638661451682
    $ws.Range("I5").Value = 0.3656254573230189
    $ws.Range("J5").Value = 0.365625457323019
    $ws.Range("M5").Value = 9.084137666666667
    $ws.Range("N5").Value = 27.252413
    $ws.Range("O5").Value = 0.2765376761551382
    $ws.Range("P5").Value = 0.2765376761551382
    $ws.Range("Q5").Value = 521.8237536414
    $ws.Range("R5").Value = 4696.4137827726
    $ws.Range("S5").Value = 0.1011092143112673
    $ws.Range("T5").Value = 0.1011092143112673
    $ws.Range("I6").Value = 0.3656254573230189
    $ws.Range("J6").Value = 0.365625457323019
    $ws.Range("O6").Value = 0.3707916163717078
    $ws.Range("P6").Value = 0.3707916163717078
    $ws.Range("R6").Value = 6297.119733832201
    $ws.Range("S6").Value = 0.1355708543074471
    $ws.Range("T6").Value = 0.1355708543074471
    $ws.Range("I7").Value = 0.3656254573230189
    $ws.Range("J7").Value = 0.365625457323019
    $ws.Range("M7").Value = 11.58507333333333
    $ws.Range("N7").Value = 34.75522
    $ws.Range("O7").Value = 0.3526707074731541
    $ws.Range("P7").Value = 0.3526707074731541
    $ws.Range("Q7").Value = 665.486001516
    $ws.Range("R7").Value = 5989.374013644
    $ws.Range("S7").Value = 0.1289453887043046
    $ws.Range("T7").Value = 0.1289453887043046
    $ws.Range("G8").Value = 35.755375
    $ws.Range("H8").Value = 107.266125
    $ws.Range("I8").Value = 0.2275818516336261
    $ws.Range("J8").Value = 0.2275818516336262
    $ws.Range("M8").Value = 9.084137666666667
    $ws.Range("N8").Value = 27.252413
    $ws.Range("O8").Value = 0.2765376761551382
    $ws.Range("P8").Value = 0.2765376761551382
    $ws.Range("Q8").Value = 324.8067488232917
    $ws.Range("R8").Value = 2923.260739409625
    $ws.Range("S8").Value = 0.06293495638584641
    $ws.Range("T8").Value = 0.06293495638584641
    $ws.Range("G9").Value = 35.755375
    $ws.Range("H9").Value = 107.266125
    $ws.Range("I9").Value = 0.2275818516336261
    $ws.Range("J9").Value = 0.2275818516336262
    $ws.Range("O9").Value = 0.3707916163717078
    $ws.Range("P9").Value = 0.3707916163717078
    $ws.Range("Q9").Value = 435.5125170613751
    $ws.Range("R9").Value = 3919.612653552375
    $ws.Range("S9").Value = 0.08438544262409843
    $ws.Range("T9").Value = 0.08438544262409843
    $ws.Range("G10").Value = 35.755375
    $ws.Range("H10").Value = 107.266125
    $ws.Range("I10").Value = 0.2275818516336261
    $ws.Range("J10").Value = 0.2275818516336262
    $ws.Range("M10").Value = 11.58507333333333
    $ws.Range("N10").Value = 34.75522
    $ws.Range("O10").Value = 0.3526707074731541
    $ws.Range("P10").Value = 0.3526707074731541
    $ws.Range("Q10").Value = 414.2286414358334
    $ws.Range("R10").Value = 3728.0577729225
    $ws.Range("S10").Value = 0.08026145262368131
    $ws.Range("T10").Value = 0.08026145262368133
